$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.195.34'
$ws.Range("E2").Value = '  +3.20%  '

$ws.Range("D3").Value = '1.581.27'
$ws.Range("E3").Value = '  +1.95%  '

$ws.Range("E4").Value = '  -0.32%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.Value = '''213.20'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("E6").Value = '  +6.74%  '

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.Value = '''0.997'
$cell.Style = $origStyle
$ws.Range("E7").Value = '  -0.37%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.Value = '''26.40'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  +11.05%  '

$ws.Range("E9").Value = '  +2.47%  '

$ws.Range("E10").Value = '  +1.90%  '

$ws.Range("E11").Value = '  +1.72%  '

$ws.Range("D12").Value = '1.806.31'
$ws.Range("E12").Value = '  +1.87%  '

$ws.Range("D13").Value = '1.582.98'
$ws.Range("E13").Value = '  +2.08%  '

$ws.Range("D14").Value = '29.238.20'
$ws.Range("E14").Value = '  +3.30%  '

$ws.Range("E15").Value = '  +2.90%  '

$ws.Range("E16").Value = '  +2.76%  '

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.Value = '''62.77'
$cell.Style = $origStyle
$ws.Range("E17").Value = '  +3.25%  '

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.Value = '''238.11'
$cell.Style = $origStyle
$ws.Range("E18").Value = '  +4.68%  '

$ws.Range("E19").Value = '  +1.83%  '

$ws.Range("E20").Value = '  +2.24%  '

$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("E22").Value = '  +1.83%  '

$ws.Range("E23").Value = '  +3.02%  '

$ws.Range("E24").Value = '  +3.23%  '

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.Value = '''154.22'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  +2.29%  '

$ws.Range("E26").Value = '  +5.16%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.Value = '''15.17'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +2.90%  '

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.Value = '''6.38'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +2.28%  '

$ws.Range("E30").Value = '  +0.45%  '

$ws.Range("E31").Value = '  +0.85%  '

$ws.Range("E32").Value = '  +1.82%  '

$ws.Range("D33").Value = '1.426.47'
$ws.Range("E33").Value = '  +2.83%  '

$ws.Range("E34").Value = '  +2.20%  '

$ws.Range("E35").Value = '  -3.17%  '

$ws.Range("E36").Value = '  +2.02%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.Value = '''2.77'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +7.44%  '

$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("E39").Value = '  +2.22%  '

$ws.Range("E40").Value = '  +3.75%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.Value = '''1.96'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  +3.00%  '

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.Value = '''54.54'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +27.89%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.Value = '''0.998'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  -0.34%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.Value = '''0.792'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +2.20%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.Value = '''0.0471'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +2.95%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.Value = '''64.61'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +4.47%  '

$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("D48").Value = '1.718.84'
$ws.Range("E48").Value = '  +2.01%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.Value = '''0.839'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  -6.39%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.Value = '''85.47'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.Value = '''0.0514'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +0.81%  '
